# Rename the "HU Transfer" sheet to "Allocations", make it the active sheet
# (with a new selection at G16), and clear the previously-active "Picking"
# sheet's tabSelected flag.

$wb = $excel.ActiveWorkbook

# Rename "HU Transfer" -> "Allocations"
$huSheet = $wb.Worksheets.Item("HU Transfer")
$huSheet.Name = "Allocations"

# Update selection on the Allocations sheet and make it the active tab.
$huSheet.Select()
$huSheet.Range("G16").Select()

$wb.Save()
